$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the results data (rows 2-25) for columns B,C,E,F,G,H,I,J,K,L
# D column remains unchanged (all zeros).

$ws.Range("B2").Value = 0.7011252215273771
$ws.Range("C2").Value = 0.07441954820658481
$ws.Range("E2").Value = 0.1672815938945575
$ws.Range("F2").Value = 2.485942218732831
$ws.Range("G2").Value = 1.405696287753514
$ws.Range("H2").Value = 1.314022443498715
$ws.Range("I2").Value = 1.353706269373127
$ws.Range("J2").Value = 0.08815242923239985
$ws.Range("K2").Value = 0.4145962733496731
$ws.Range("L2").Value = 0.4185993786323934

$ws.Range("B3").Value = 0.6633612166925786
$ws.Range("C3").Value = 0.07235220454283109
$ws.Range("E3").Value = 0.1663215372298623
$ws.Range("F3").Value = 2.483043847407188
$ws.Range("G3").Value = 1.409841234333612
$ws.Range("H3").Value = 1.320986855299353
$ws.Range("I3").Value = 1.36117835806521
$ws.Range("J3").Value = 0.08808919817806959
$ws.Range("K3").Value = 0.3800044439820454
$ws.Range("L3").Value = 0.4103376167998221

$ws.Range("B4").Value = 0.640451261161985
$ws.Range("C4").Value = 0.07106141399443544
$ws.Range("E4").Value = 0.165788643321708
$ws.Range("F4").Value = 2.482440281539851
$ws.Range("G4").Value = 1.413137928468402
$ws.Range("H4").Value = 1.32578460037449
$ws.Range("I4").Value = 1.366359558901493
$ws.Range("J4").Value = 0.08805538167351123
$ws.Range("K4").Value = 0.3588943865144074
$ws.Range("L4").Value = 0.4054411284788131

$ws.Range("B5").Value = 0.63118556011392
$ws.Range("C5").Value = 0.07052999779362068
$ws.Range("E5").Value = 0.165585765772903
$ws.Range("F5").Value = 2.482490223302221
$ws.Range("G5").Value = 1.414670236810522
$ws.Range("H5").Value = 1.327870919343667
$ws.Range("I5").Value = 1.368620125985458
$ws.Range("J5").Value = 0.08804287181959758
$ws.Range("K5").Value = 0.3503247850771771
$ws.Range("L5").Value = 0.4034902696019458

$ws.Range("B6").Value = 0.6296512596435093
$ws.Range("C6").Value = 0.07044142932247865
$ws.Range("E6").Value = 0.1655529421931625
$ws.Range("F6").Value = 2.482516392846264
$ws.Range("G6").Value = 1.41493607924248
$ws.Range("H6").Value = 1.328225276487558
$ws.Range("I6").Value = 1.369004501698917
$ws.Range("J6").Value = 0.08804087164989127
$ws.Range("K6").Value = 0.3489038082497586
$ws.Range("L6").Value = 0.4031690238547014

$ws.Range("B7").Value = 0.6403260153180099
$ws.Range("C7").Value = 0.07105426905290813
$ws.Range("E7").Value = 0.1657858493507582
$ws.Range("F7").Value = 2.482439756758552
$ws.Range("G7").Value = 1.413157829137873
$ws.Range("H7").Value = 1.325812205924791
$ws.Range("I7").Value = 1.36638944166458
$ws.Range("J7").Value = 0.08805520780139631
$ws.Range("K7").Value = 0.3587786801028869
$ws.Range("L7").Value = 0.405414638100055

$ws.Range("B8").Value = 0.6880469636912778
$ws.Range("C8").Value = 0.07371116461191463
$ws.Range("E8").Value = 0.1669388529417262
$ws.Range("F8").Value = 2.484698869016142
$ws.Range("G8").Value = 1.406969369081651
$ws.Range("H8").Value = 1.316315548230349
$ws.Range("I8").Value = 1.356159499624681
$ws.Range("J8").Value = 0.08812959460892245
$ws.Range("K8").Value = 0.4026423457349608
$ws.Range("L8").Value = 0.4157142409710843

$ws.Range("B9").Value = 0.7838085046182357
$ws.Range("C9").Value = 0.0787523199024065
$ws.Range("E9").Value = 0.1696468170010874
$ws.Range("F9").Value = 2.498455972348324
$ws.Range("G9").Value = 1.400805481798741
$ws.Range("H9").Value = 1.301829441374082
$ws.Range("I9").Value = 1.340807085547063
$ws.Range("J9").Value = 0.08831472109466532
$ws.Range("K9").Value = 0.4896748508873259
$ws.Range("L9").Value = 0.4373040408796811

$ws.Range("B10").Value = 0.8554761184971653
$ws.Range("C10").Value = 0.08235476449375767
$ws.Range("E10").Value = 0.1719064377132575
$ws.Range("F10").Value = 2.514248150151445
$ws.Range("G10").Value = 1.399929399386863
$ws.Range("H10").Value = 1.293706915078317
$ws.Range("I10").Value = 1.332400028904807
$ws.Range("J10").Value = 0.08847405440200262
$ws.Range("K10").Value = 0.5542288047392674
$ws.Range("L10").Value = 0.4540085448982012

$ws.Range("B11").Value = 0.8883611272579515
$ws.Range("C11").Value = 0.08397200617872613
$ws.Range("E11").Value = 0.1729925645116879
$ws.Range("F11").Value = 2.522666917606116
$ws.Range("G11").Value = 1.400326672038815
$ws.Range("H11").Value = 1.290558805313566
$ws.Range("I11").Value = 1.329199640730735
$ws.Range("J11").Value = 0.08855147803896557
$ws.Range("K11").Value = 0.5837275131464423
$ws.Range("L11").Value = 0.4617895983486449

$ws.Range("B12").Value = 0.9008540840765136
$ws.Range("C12").Value = 0.08458134269911
$ws.Range("E12").Value = 0.1734121764914889
$ws.Range("F12").Value = 2.526032368167492
$ws.Range("G12").Value = 1.400591735698498
$ws.Range("H12").Value = 1.289445312540352
$ws.Range("I12").Value = 1.328077504239801
$ws.Range("J12").Value = 0.08858149629843126
$ws.Range("K12").Value = 0.594916744284717
$ws.Range("L12").Value = 0.46476212054732

$ws.Range("B13").Value = 0.8981617273237248
$ws.Range("C13").Value = 0.08445024777754639
$ws.Range("E13").Value = 0.1733214363617428
$ws.Range("F13").Value = 2.52529966813475
$ws.Range("G13").Value = 1.400529548505304
$ws.Range("H13").Value = 1.289681626293159
$ws.Range("I13").Value = 1.328315182620045
$ws.Range("J13").Value = 0.08857500038886101
$ws.Range("K13").Value = 0.5925061171001289
$ws.Range("L13").Value = 0.4641207809096386

$ws.Range("B14").Value = 0.8893881296521329
$ws.Range("C14").Value = 0.08402219824155566
$ws.Range("E14").Value = 0.1730269198058174
$ws.Range("F14").Value = 2.522940239657984
$ws.Range("G14").Value = 1.400346180645627
$ws.Range("H14").Value = 1.290465621829696
$ws.Range("I14").Value = 1.329105522269451
$ws.Range("J14").Value = 0.08855393370084741
$ws.Range("K14").Value = 0.5846476851341436
$ws.Range("L14").Value = 0.4620336293564975

$ws.Range("B15").Value = 0.8840192549633628
$ws.Range("C15").Value = 0.08375960525510351
$ws.Range("E15").Value = 0.1728476018009317
$ws.Range("F15").Value = 2.521518126718462
$ws.Range("G15").Value = 1.400248795717559
$ws.Range("H15").Value = 1.290956081116036
$ws.Range("I15").Value = 1.329601321671952
$ws.Range("J15").Value = 0.08854112053070651
$ws.Range("K15").Value = 0.5798365939340329
$ws.Range("L15").Value = 0.4607585702496095

$ws.Range("B16").Value = 0.8533326514832709
$ws.Range("C16").Value = 0.08224864227196349
$ws.Range("E16").Value = 0.1718366237575495
$ws.Range("F16").Value = 2.513722801470934
$ws.Range("G16").Value = 1.399919466523286
$ws.Range("H16").Value = 1.29392365489403
$ws.Range("I16").Value = 1.332621742782955
$ws.Range("J16").Value = 0.0884690930845089
$ws.Range("K16").Value = 0.5523036339126577
$ws.Range("L16").Value = 0.4535036836304016

$ws.Range("B17").Value = 0.8345794932134254
$ws.Range("C17").Value = 0.08131621422177915
$ws.Range("E17").Value = 0.1712312945812826
$ws.Range("F17").Value = 2.509256804978293
$ws.Range("G17").Value = 1.399921393602469
$ws.Range("H17").Value = 1.295884223385116
$ws.Range("I17").Value = 1.334634534312592
$ws.Range("J17").Value = 0.08842616445905804
$ws.Range("K17").Value = 0.5354467910526068
$ws.Range("L17").Value = 0.4490995564128468

$ws.Range("B18").Value = 0.8238198494631774
$ws.Range("C18").Value = 0.08077787883096477
$ws.Range("E18").Value = 0.1708886077110492
$ws.Range("F18").Value = 2.506804316849696
$ws.Range("G18").Value = 1.399997395827938
$ws.Range("H18").Value = 1.297063364945387
$ws.Range("I18").Value = 1.335850972015862
$ws.Range("J18").Value = 0.0884019388755668
$ws.Range("K18").Value = 0.52576370262193
$ws.Range("L18").Value = 0.4465835687891513

$ws.Range("B19").Value = 0.8201814186282093
$ws.Range("C19").Value = 0.0805952589587946
$ws.Range("E19").Value = 0.170773523023243
$ws.Range("F19").Value = 2.50599391364392
$ws.Range("G19").Value = 1.400035986201488
$ws.Range("H19").Value = 1.297471443147202
$ws.Range("I19").Value = 1.336272922985657
$ws.Range("J19").Value = 0.08839381685537262
$ws.Range("K19").Value = 0.5224873376468793
$ws.Range("L19").Value = 0.445734649559526

$ws.Range("B20").Value = 0.8365730422169406
$ws.Range("C20").Value = 0.08141568248228737
$ws.Range("E20").Value = 0.1712951658671855
$ws.Range("F20").Value = 2.509720188918195
$ws.Range("G20").Value = 1.399913435857826
$ws.Range("H20").Value = 1.295670190190208
$ws.Range("I20").Value = 1.334414190202615
$ws.Range("J20").Value = 0.08843068615303906
$ws.Range("K20").Value = 0.5372399375958423
$ws.Range("L20").Value = 0.4495666096214279

$ws.Range("B21").Value = 0.891964063821888
$ws.Range("C21").Value = 0.08414801015349127
$ws.Range("E21").Value = 0.1731132010848633
$ws.Range("F21").Value = 2.523628446122331
$ws.Range("G21").Value = 1.400396927832546
$ws.Range("H21").Value = 1.290233209476881
$ws.Range("I21").Value = 1.328870943347042
$ws.Range("J21").Value = 0.0885601025980165
$ws.Range("K21").Value = 0.5869553923490116
$ws.Range("L21").Value = 0.4626459719795548

$ws.Range("B22").Value = 0.9283988664636809
$ws.Range("C22").Value = 0.08591581282468752
$ws.Range("E22").Value = 0.1743498597456608
$ws.Range("F22").Value = 2.533752531404048
$ws.Range("G22").Value = 1.401381128296123
$ws.Range("H22").Value = 1.287138128270669
$ws.Range("I22").Value = 1.325771445839955
$ws.Range("J22").Value = 0.08864875654316862
$ws.Range("K22").Value = 0.6195561984125106
$ws.Range("L22").Value = 0.4713455927169292

$ws.Range("B23").Value = 0.9089317473348615
$ws.Range("C23").Value = 0.08497393803735065
$ws.Range("E23").Value = 0.1736854140165391
$ws.Range("F23").Value = 2.528254523861392
$ws.Range("G23").Value = 1.400794638292965
$ws.Range("H23").Value = 1.288748100316994
$ws.Range("I23").Value = 1.32737780518481
$ws.Range("J23").Value = 0.08860107122082539
$ws.Range("K23").Value = 0.6021466982909374
$ws.Range("L23").Value = 0.4666886377585371

$ws.Range("B24").Value = 0.8356716907746033
$ws.Range("C24").Value = 0.0813707199542506
$ws.Range("E24").Value = 0.1712662730689303
$ws.Range("F24").Value = 2.509510334581464
$ws.Range("G24").Value = 1.399916800271924
$ws.Range("H24").Value = 1.295766792661425
$ws.Range("I24").Value = 1.334513623190745
$ws.Range("J24").Value = 0.08842864047733201
$ws.Range("K24").Value = 0.5364292306766458
$ws.Range("L24").Value = 0.4493554050044821

$ws.Range("B25").Value = 0.7576708583632126
$ws.Range("C25").Value = 0.0774064678830868
$ws.Range("E25").Value = 0.1688666410278437
$ws.Range("F25").Value = 2.493735911503734
$ws.Range("G25").Value = 1.401832423859545
$ws.Range("H25").Value = 1.305305600277279
$ws.Range("I25").Value = 1.344456045636775
$ws.Range("J25").Value = 0.08826050058758561
$ws.Range("K25").Value = 0.4660223430799988
$ws.Range("L25").Value = 0.4313150268037731
